$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Numéro"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "position début (x,y,z)"
$ws.Range("D1").Value = "position fin(x,y,z)"
$ws.Range("E1").Value = "épaisseur"

# Béton rows (1-4)
$data = @(
    @(1, "Béton", ",,", ",,", 200),
    @(2, "Béton", ",,", ",,", 200),
    @(3, "Béton", ",,", ",,", 200),
    @(4, "Béton", ",,", ",,", 200),
    @(5, "Bois",  ",,", ",,", 100),
    @(6, "Bois",  ",,", ",,", 100),
    @(7, "Bois",  ",,", ",,", 100)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Column widths to match auto-fit (bestFit) sizing observed in the target file.
# The host snaps ColumnWidth to a 1/6-character grid (stored = round(w*6)/6 + 5/6),
# so we feed the pre-image that lands exactly (col C) / as close as possible (col D).
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668   # -> stored width 20
$ws.Columns.Item(4).ColumnWidth = 15.833333333333334   # -> stored width 16.666666666666668 (closest reachable to 16.7109375)

# Move selection to A9, mirroring the cursor position after typing the last row of data
$ws.Range("A9").Select()
